$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.719.49'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '3.165.64'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.04%  '
$c = $ws.Range('D5')
$c.Value = "'612.27"
$c.ClearFormats()
$ws.Range('E5').Value = '  +1.39%  '
$c = $ws.Range('D6')
$c.Value = "'145.90"
$c.ClearFormats()
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('D8').Value = '3.158.96'
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('E10').Value = '  +0.99%  '
$c = $ws.Range('D11')
$c.Value = "'5.45"
$c.ClearFormats()
$ws.Range('E11').Value = '  -1.41%  '
$c = $ws.Range('D12')
$c.Value = "'0.473"
$c.ClearFormats()
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('E13').Value = '  +1.38%  '
$c = $ws.Range('D14')
$c.Value = "'35.64"
$c.ClearFormats()
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').Value = '3.686.93'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('E16').Value = '  +3.28%  '
$ws.Range('D17').Value = '64.677.19'
$ws.Range('E17').Value = '  +0.80%  '
$ws.Range('D18').Value = '3.163.98'
$ws.Range('E18').Value = '  +0.49%  '
$c = $ws.Range('D19')
$c.Value = "'6.86"
$c.ClearFormats()
$ws.Range('E19').Value = '  -0.64%  '
$c = $ws.Range('D20')
$c.Value = "'479.95"
$c.ClearFormats()
$ws.Range('E20').Value = '  +0.39%  '
$c = $ws.Range('D21')
$c.Value = "'14.60"
$c.ClearFormats()
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('E22').Value = '  +2.49%  '
$c = $ws.Range('D23')
$c.Value = "'7.91"
$c.ClearFormats()
$ws.Range('E23').Value = '  +3.31%  '
$c = $ws.Range('D24')
$c.Value = "'13.74"
$c.ClearFormats()
$ws.Range('E24').Value = '  +0.52%  '
$c = $ws.Range('D25')
$c.Value = "'84.08"
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.92%  '
$c = $ws.Range('D26')
$c.Value = "'1.00"
$c.ClearFormats()
$ws.Range('E26').Value = '  +0.03%  '
$c = $ws.Range('D27')
$c.Value = "'8.72"
$c.ClearFormats()
$ws.Range('E27').Value = '  +3.28%  '
$ws.Range('E28').Value = '  -3.66%  '
$c = $ws.Range('D29')
$c.Value = "'7.12"
$c.ClearFormats()
$ws.Range('E29').Value = '  +4.32%  '
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('E31').Value = '  -5.08%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('E33').Value = '  -0.95%  '
$c = $ws.Range('D34')
$c.Value = "'26.54"
$c.ClearFormats()
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('E35').Value = '  +1.66%  '
$ws.Range('D36').Value = '0.0₃0791'
$ws.Range('E36').Value = '  +7.82%  '
$ws.Range('E37').Value = '  -0.70%  '
$c = $ws.Range('D38')
$c.Value = "'53.27"
$c.ClearFormats()
$ws.Range('E38').Value = '  -2.07%  '
$c = $ws.Range('D39')
$c.Value = "'3.18"
$c.ClearFormats()
$ws.Range('E39').Value = '  +2.63%  '
$c = $ws.Range('D40')
$c.Value = "'460.90"
$c.ClearFormats()
$ws.Range('E40').Value = '  +2.45%  '
$ws.Range('E41').Value = '  +0.70%  '
$c = $ws.Range('D42')
$c.Value = "'0.119"
$c.ClearFormats()
$ws.Range('E42').Value = '  -2.31%  '
$c = $ws.Range('D43')
$c.Value = "'8.32"
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.64%  '
$ws.Range('D44').Value = '2.856.53'
$ws.Range('E44').Value = '  +0.06%  '
$c = $ws.Range('D45')
$c.Value = "'2.32"
$c.ClearFormats()
$ws.Range('E45').Value = '  +3.02%  '
$c = $ws.Range('D46')
$c.Value = "'0.267"
$c.ClearFormats()
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('E47').Value = '  +6.38%  '
$c = $ws.Range('D48')
$c.Value = "'26.58"
$c.ClearFormats()
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('E49').Value = '  +0.05%  '
$c = $ws.Range('D50')
$c.Value = "'35.76"
$c.ClearFormats()
$ws.Range('E50').Value = '  +7.81%  '
$ws.Range('E51').Value = '  -0.01%  '
